# 1st May 2025 - Final Changes
# Add a new "Warning" sheet (right after "Tooltip") that carries the same
# look & feel as the existing "Tooltip" sheet (bold header in A1, wrapped
# message in A2), then make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$tooltip = $wb.Worksheets.Item("Tooltip")

# Duplicate the Tooltip sheet right after itself - this carries over the
# same column width / wrap-text / bold styling as a starting point so we
# only need to touch the bits that actually differ.
$tooltip.Copy([System.Reflection.Missing]::Value, $tooltip)
$warning = $wb.Worksheets.Item("Tooltip (2)")
$warning.Name = "Warning"

# Replace the header + body text with the new warning copy.
$warning.Range("A1").Value = "Message"
$warning.Range("A2").Value = "An asset is typically considered a potential round trip if it is an operating company acquired either by a private equity firm or by a PE-owned operating company. This company is not listed as an Operating Company that is PE-owned. If you still want to consider them a round trip candidate no change is needed; otherwise, please change the selection."

# Widen column A to fit the longer message and grow row 2 to fit the wrapped text.
$warning.Columns("A:A").ColumnWidth = 96.25
$warning.Rows(2).RowHeight = 43.2

# Leave the selection parked on A8, matching the saved view state.
$warning.Range("A8").Select() | Out-Null

# The newly added/copied sheet becomes the active tab.
$warning.Activate() | Out-Null
